# household_member.xlsx edit:
# - Remove the " for {{instance_name}}" suffix from the note prompt on the
#   survey sheet (row 6, column D) and add a new 'instance_name' setting
#   row to the settings sheet that points at the 'member_name' field.
# - Active sheet/tab moves from 'survey' to 'settings', with the selection
#   on 'survey' moving from E7 to D7, and the selection on 'settings'
#   moving from C9 to B18.
# - The survey sheet's note row (row 6) grows taller to fit the now
#   multi-line-looking prompt text.

$wb = $excel.ActiveWorkbook

$surveySheet = $wb.Worksheets.Item("survey")
$settingsSheet = $wb.Worksheets.Item("settings")

# Update the note prompt text (drop " for {{instance_name}}").
$surveySheet.Range("D6").Value = "{{member_name}} age is {{evaluate calculates.ageIsOddOrEven}} in {{setting 'table_id'}}"

# Grow row 6 to its new height.
$surveySheet.Rows.Item(6).RowHeight = 62.5

# Add the new 'instance_name' setting row, referencing 'member_name'.
$settingsSheet.Range("A6").Value = "instance_name"
$settingsSheet.Range("B6").Value = "member_name"

# Update selections: survey E7 -> D7, settings C9 -> B18.
$surveySheet.Range("D7").Select() | Out-Null
$settingsSheet.Range("B18").Select() | Out-Null

# Move the active tab/selection to the settings sheet.
$settingsSheet.Activate() | Out-Null
